# Commit: "banging my head against the wall with tidyeval"
# Adds a "Sheet2" lookup/summary sheet that maps a smaller set of EDUCD
# codes to their Sheet1 Label/Group via VLOOKUP, restyles the code column,
# and sets the new sheet as the active tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- New sheet, placed immediately after Sheet1 --------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Code"
$ws2.Range("B1").Value = "IPUMS CAT"
$ws2.Range("C1").Value = "My Category"

# Left-align the header's code cell (style later shared with the plain
# code cells below it).
$ws2.Range("A1").HorizontalAlignment = -4131

# The distinct EDUCD codes this summary sheet reports on.
$codes = @(1, 2, 11, 12, 14, 15, 16, 17, 22, 23, 25, 26, 30, 40, 50, 61, 63, 64, 65, 71, 81, 101, 114, 115, 116)

$row = 2
foreach ($code in $codes) {
    $ws2.Cells.Item($row, 1).Value = $code
    $ws2.Cells.Item($row, 2).Formula = "=VLOOKUP(A$row, Sheet1!`$A`$2:`$C`$45, 2, FALSE)"
    $ws2.Cells.Item($row, 3).Formula = "=VLOOKUP(A$row, Sheet1!`$A`$2:`$C`$45, 3, FALSE)"
    $ws2.Cells.Item($row, 1).HorizontalAlignment = -4131
    $row = $row + 1
}

# First code cell (A2) gets the distinct "Lucida Console" style.
$codeFont = $ws2.Range("A2").Font
$codeFont.Name = "Lucida Console"
$codeFont.Size = 10
$codeFont.Color = 0
$codeFont.Family = 3
$ws2.Range("A2").VerticalAlignment = -4108

$ws2.Range("A1:C26").Select()

# --- Sheet1 cosmetic tweaks ------------------------------------------------
$sheet1.PageSetup.Orientation = 1

# --- Workbook-level: make the new sheet the active tab --------------------
$ws2.Activate()
